$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns C..H (processing, completed, canceled, deferred, closed, new_or_reopened)
# in the two data rows (2 and 3) get a ":formatN()" formatter suffix appended to
# their template tag, and their number format switched from Text to a plain
# integer number format so the templating engine's numeric output renders
# as a number instead of text.
$cols = @("C", "D", "E", "F", "G", "H")
$fields = @("processing", "completed", "canceled", "deferred", "closed", "new_or_reopened")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $col = $cols[$i]
    $field = $fields[$i]

    $cellRow2 = $col + "2"
    $cellRow3 = $col + "3"

    $ws.Range($cellRow2).Value = "{d.tickets[i]." + $field + ":formatN()}"
    $ws.Range($cellRow3).Value = "{d.tickets[i+1]." + $field + ":formatN()}"

    $ws.Range($cellRow2).NumberFormat = "0"
    $ws.Range($cellRow3).NumberFormat = "0"
}
